# Rename model "Region" to "Zone" for clarity.
$wb = $excel.ActiveWorkbook

# Rename the "Region" worksheet to "Zone".
$ws = $wb.Worksheets.Item("Region")
$ws.Name = "Zone"

# Bring the renamed sheet to the foreground (this is the sheet the author
# was working on when performing the rename), and move the cell selection
# to D3.
$ws.Activate()
$ws.Range("D3").Select()
